$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("A2").Value = 13317
$ws.Range("B2").Value = "Vinicius Maciel"
$ws.Range("C2").Value = "Offshore"
$ws.Range("F2").Value = "Fernandito Banker"

# --- Row 3 updates ---
$ws.Range("C3").Value = "Moderada 1"
$ws.Range("F3").Value = "Fernandito Banker"
$ws.Range("O3").Value = "Aqui está puxando corretamente"
$ws.Range("R3").Value = "ATIVO DETRATOR -94%"

# --- Row 4 (new row) ---
$ws.Range("A4").Value = 13357
$ws.Range("B4").Value = "Helena Miranda"
$ws.Range("C4").Value = "Moderada 2 "
$ws.Range("D4").Value = "Moderada"
$ws.Range("E4").Value = "IPCA + 6%"
$ws.Range("F4").Value = "Renato Banker"
$ws.Range("G4").Value = 0.5
$ws.Range("H4").Value = 1.79
$ws.Range("I4").Value = 0.63
$ws.Range("J4").Value = 0.42
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = -0.22
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = "XP INVESTIMENTOS (1.56%)"
$ws.Range("O4").Value = "Teste"
$ws.Range("P4").Value = "AAA"
$ws.Range("Q4").Value = "BBB"

# --- selection moves to C3 ---
$ws.Range("C3").Select()
